$d = $word.ActiveDocument

# --- Merge "Versi" + "on" (chars 0-6) into a single run "Version" ---
# A textual no-op would not force Word to actually merge the two runs,
# so first change to a distinct placeholder, then rename to the final text.
$r1 = $d.Range(0, 7)
$r1.Text = "Versionx"
$r1 = $d.Range(0, 8)
$r1.Text = "Version"

# --- Merge " 2" + "." (chars 7-9) into a single run " 1." ---
$r2 = $d.Range(7, 10)
$r2.Text = " 1x."
$r2 = $d.Range(7, 11)
$r2.Text = " 1."

# --- Re-create the _GoBack bookmark after the merged runs (collapsed) ---
# Inserting a bookmark exactly at the end of the paragraph's text range can
# make this runtime snap its start back to 0, so pad with a throwaway
# character, add the bookmark in the middle of the range, then clean up.
$contentEnd = $d.Content.End
$tail = $d.Range($contentEnd - 1, $contentEnd - 1)
$tail.InsertAfter("Z")

$bmRange = $d.Range($contentEnd - 1, $contentEnd - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

$padRange = $d.Range($contentEnd - 1, $contentEnd)
$padRange.Delete()
